$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: AD = 30, AE = 31, AF = 32
$colWins = 30
$colLosses = 31
$colTies = 32

# Header row (row 1) - new headers, using the same bold/centered style (s="1")
# as the other header cells. Copy style from an existing header cell (AC1).
$ws.Cells.Item(1, $colWins).Value = "Wins"
$ws.Cells.Item(1, $colLosses).Value = "Losses"
$ws.Cells.Item(1, $colTies).Value = "Ties"

$headerStyleRange = $ws.Range("AC1")
$newHeaderRange = $ws.Range("AD1:AF1")
$newHeaderRange.Font.Bold = $headerStyleRange.Font.Bold
$newHeaderRange.HorizontalAlignment = $headerStyleRange.HorizontalAlignment
$newHeaderRange.VerticalAlignment = $headerStyleRange.VerticalAlignment
$newHeaderRange.Borders.LineStyle = $headerStyleRange.Borders.LineStyle

# Data rows 2-55: season record values (Wins=70, Losses=92, Ties=0)
$lastRow = 55
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $colWins).Value = 70
    $ws.Cells.Item($r, $colLosses).Value = 92
    $ws.Cells.Item($r, $colTies).Value = 0
}
